$d = $word.ActiveDocument

# 1. Update the title heading text (appears twice - in the title area and
#    again as the repeated "firstheader" heading). Wrap=1 (wdFindContinue)
#    plus Replace=2 (wdReplaceAll) makes Find.Execute replace every match
#    in the story in one call.
$d.Content.Find.Execute(
    "Análisis Univariado de Series de Tiempo: Identificación de Modelos y Proyecciones",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Análisis univariado de series temporales", 2) | Out-Null

# 2. Insert the new "abstract pending" paragraph right before the
#    "Palabras clave" paragraph, and restyle that surrounding pair of
#    paragraphs (new one = AbstractFirstParagraph, keywords one = Body Text).
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Palabras clave*") {
        $kwRange = $para.Range
        $kwRange.Collapse(1)
        $kwRange.InsertBefore("Este abstract será actualizado una vez que se complete el contenido final del artículo.`r")
        break
    }
}

foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Este abstract será actualizado*") {
        $para.Style = "AbstractFirstParagraph"
    }
    elseif ($para.Range.Text -like "*Palabras clave*") {
        $para.Style = "Body Text"
    }
}

# 3. Expand the "holaa" sentence into "holaaa, cómo estas por favor".
$d.Content.Find.Execute(
    "todas las secciones serán ampliadas y refinadas en futuras revisiones. holaa",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "todas las secciones serán ampliadas y refinadas en futuras revisiones. holaaa, cómo estas por favor", 2) | Out-Null

# 4. Remove the "Por Editar" bullet entry entirely (paragraph + its
#    hyperlinks) from the "Publicaciones Similares" list.
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Por Editar*") {
        $para.Range.Delete()
        break
    }
}
